# Add a new state q4 to the DFA transition table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 5 (the old blank separator row and
# everything below it shifts down by one row).
$ws.Rows(5).Insert()

# Give the new A5 header cell the same formatting (border / bold /
# centered) as the other state-header cells in column A (A2:A4).
$ws.Range("A2").Copy()
$ws.Range("A5").PasteSpecial(-4122)

# Update the transition table with the new state q4's transitions,
# and fix up the other transitions that changed.
$ws.Range("C2").Value = "q1"
$ws.Range("B3").Value = "q2"
$ws.Range("C3").Value = "q3"
$ws.Range("B4").Value = "q2"
$ws.Range("C4").Value = "q4"
$ws.Range("A5").Value = "q4"
$ws.Range("B5").Value = "q2"
$ws.Range("C5").Value = "q1"

# Update the summary rows (now shifted down to rows 7-10).
$ws.Range("B7").Value = "{q1,q2,q3,q4}"
$ws.Range("B10").Value = "{q4}"
